$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates (America De Cali vs Millonarios)
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("U4").Value = 2.5
$ws.Range("V4").Value = 1.5
$ws.Range("Y4").Value = 10
$ws.Range("AC4").Value = 5.5
$ws.Range("AG4").Value = 9

# Row 7 updates (Tepatitlan de Morelos vs Tampico Madero)
$ws.Range("G7").Value = 2.22
$ws.Range("H7").Value = 2.9
$ws.Range("I7").Value = 3.4
$ws.Range("J7").Value = 2.77
$ws.Range("L7").Value = 3.95
$ws.Range("N7").Value = 6.85
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 2.55
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.55
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.42
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 6.7
$ws.Range("X7").Value = 10.5
$ws.Range("Y7").Value = 8.75
$ws.Range("Z7").Value = 23
$ws.Range("AA7").Value = 19.5
$ws.Range("AB7").Value = 32
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 5.7
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 8.25
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 12
$ws.Range("AK7").Value = 37
$ws.Range("AL7").Value = 45
$ws.Range("AO7").Value = 11.25
$ws.Range("AP7").Value = 19
$ws.Range("AQ7").Value = 45
$ws.Range("AR7").Value = 75
$ws.Range("AT7").Value = 2.4
$ws.Range("AX7").Value = 20
$ws.Range("AY7").Value = 27
$ws.Range("AZ7").Value = 110
$ws.Range("BA7").Value = 150
$ws.Range("BB7").Value = 350

# Row 8 updates (Tacuary vs Sp. Luqueno)
$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67

# Row 13 updates (La Guaira vs Rayo Zuliano)
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 7

